# Edit: update lfo1ratesyncoff description cell (D191) with new "Edit: Benny Rönnhager..." note,
# narrow column D slightly, and update the custom row heights that reflow as a result
# (matches the target OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new paragraph (as a literal line break, like Alt+Enter in Excel) to the
# existing lfo1ratesyncoff NRPN notes cell.
$ws.Range("D191").Value = "[0,8192] seemingly only output in increments of 8, and displayed as [0.02 Hz ...150.00 Hz]. To display: if 8192, display 150.00Hz.  Else divide by 6.4 (cutting into 1280 even pieces).  Now we need to map to an exponential function to get the Hz value.  It seems the following function is a pretty close fit:  `n`n2^(1 + 0.012571 *  v) / 100  I would then display as x.xx, perhaps rounded down.  Would be nice to know what their exact function is.`nEdit: Benny Rönnhager manually went through the entire list, and reports the following values for all elements [0...8192] in multiples of 8, that is, [0...1024]:   LFO_RATES_SYNC_OFF    However I have not verified it yet."

# Column D got a bit narrower, which reflows the wrapped text in several rows below.
$ws.Columns.Item(4).ColumnWidth = 69.15

# Update the custom row heights to their new values (Excel would normally recompute
# these via wrap-text autofit once the column width / text changes; we pin them to the
# exact target heights here).
$ws.Rows.Item(51).RowHeight = 260.35
$ws.Rows.Item(127).RowHeight = 68.35
$ws.Rows.Item(167).RowHeight = 44.35
$ws.Rows.Item(184).RowHeight = 68.35
$ws.Rows.Item(191).RowHeight = 164.35
$ws.Rows.Item(193).RowHeight = 116.35
$ws.Rows.Item(569).RowHeight = 200.35
$ws.Rows.Item(572).RowHeight = 212.35
$ws.Rows.Item(574).RowHeight = 212.35
$ws.Rows.Item(692).RowHeight = 56.35
$ws.Rows.Item(1123).RowHeight = 116.35
$ws.Rows.Item(1129).RowHeight = 248.35
$ws.Rows.Item(1180).RowHeight = 56.35
$ws.Rows.Item(1181).RowHeight = 32.35
$ws.Rows.Item(1182).RowHeight = 32.35
$ws.Rows.Item(1183).RowHeight = 32.35
$ws.Rows.Item(1184).RowHeight = 32.35
$ws.Rows.Item(1185).RowHeight = 56.35
$ws.Rows.Item(1186).RowHeight = 32.35
$ws.Rows.Item(1187).RowHeight = 32.35
$ws.Rows.Item(1188).RowHeight = 32.35
$ws.Rows.Item(1189).RowHeight = 32.35
$ws.Rows.Item(1309).RowHeight = 56.35
